$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The C column (meeting place) style for filled-in rows uses a thicker/different
# border than the blank template rows; copy that formatting from the row above
# (row 9, an already-filled attendance row) onto C11 before writing the value.
$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial(-4122)

# Fill in row 11 with the 9/8 meeting's attendance data
$ws.Range("B11").Value = "9/8/ 4:15"
$ws.Range("C11").Value = "Google Hangout"
$ws.Range("D11").Value = "A"
$ws.Range("E11").Value = "A"
$ws.Range("F11").Value = "A"
$ws.Range("G11").Value = "A"
$ws.Range("H11").Value = "A"
$ws.Range("I11").Value = "A"

# Update the selected cell to I11 as per the saved view state
$ws.Range("I11").Select()
